# Add a new column K ("obs" header / "Apa" values) to Sheet1 -- mirrors the
# upstream upload, which appended a 12th column pulling in two more shared
# strings, a header cell styled like the rest of row 1, and a distinctly
# styled run of data cells below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1: give it the same (bold/centered/wrapped) look as the other
# header cells by copying J1's format instead of reassembling it by hand --
# this reuses the existing header style exactly, with no stray new styles.
$ws.Range("K1").Value = "obs"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $False

# Data cells K2:K6: new text value "Apa", given its own explicit font so the
# workbook gets a genuinely new font + cell style (instead of silently
# reusing the default "no style" formatting of a plain cell).
$ws.Range("K2:K6").Value = "Apa"
$ws.Range("K2:K6").Font.Name = "Arial"
$ws.Range("K2:K6").Font.Size = 9
$ws.Range("K2:K6").Font.Bold = $False
$ws.Range("K2:K6").Font.Italic = $False

# Move the active selection the way the recorded workbook shows (I10 -> I12).
$ws.Range("I12").Select() | Out-Null
